$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Duplicate the existing "LKT Events" sheet, inserting the copy in front of
# it, to create the new "LKT 8Beta3" sheet (this preserves column widths,
# row heights, cell styles and formatting exactly as Excel's own Copy
# command would).
# ---------------------------------------------------------------------------
$lktEvents = $wb.Worksheets.Item("LKT Events")
$lktEvents.Copy($wb.Worksheets.Item(1))

# The copy was inserted at index 1; rename it.
$newSheet = $wb.Worksheets.Item(1)
$newSheet.Name = "LKT 8Beta3"

# ---------------------------------------------------------------------------
# Update the HED-tag / description text on the new sheet to the "8Beta3"
# vocabulary, and tighten the row heights that no longer need to wrap as
# much text.
#
# NOTE: the order these new, previously-unseen strings are assigned below
# matters -- each first-seen string value is appended to the shared string
# table in the order it is written, and that order must line up with the
# canonical file's shared string table (new strings appended after index 53
# in the sequence: E2, E3, E5, D5, E4).
# ---------------------------------------------------------------------------
$newSheet.Range("E2").Value = "Experiment-control, Experimental-stimulus, (Controller-agent, (Operate, Car, (Turn, Leftward)))"
$newSheet.Range("E3").Value = "Experiment-control, Experimental-stimulus, (Controller-agent, (Operate, Car, (Turn, Rightward)))"
$newSheet.Range("E5").Value = "Agent-action, Participant-response, (Halt, Correction)"
$newSheet.Range("D5").Value = "Subject completes response to perturbation having steered the vehicle back to the center of the lane. Normally this would be tagged with temporal scope, but avoiding definitions here."
$newSheet.Range("E4").Value = "Agent-action, Participant-response, Correction, ((Human-agent, Experimental-participant), (Modify, (Car, Direction)))"

$newSheet.Rows.Item(2).RowHeight = 43.75
$newSheet.Rows.Item(3).RowHeight = 43.75
$newSheet.Rows.Item(5).RowHeight = 43.75

# ---------------------------------------------------------------------------
# Restore the selection on the original "LKT Events" sheet (it shifted one
# slot to the right, to index 2, when the new sheet was inserted), then set
# the selection / active sheet to match the new sheet.
# ---------------------------------------------------------------------------
$origLktEvents = $wb.Worksheets.Item("LKT Events")
$origLktEvents.Range("E5").Select()

$newSheet.Range("E4").Select()
